$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.146.39"
$ws.Range("E2").Value = "  +3.74%  "
$ws.Range("D3").Value = "'3.230.47"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'576.22"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("D6").Value = "'180.66"
$ws.Range("E6").Value = "  +5.68%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "  -4.12%  "
$ws.Range("D9").Value = "'3.229.82"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("E10").Value = "  +3.48%  "
$ws.Range("E11").Value = "  +3.20%  "
$ws.Range("E12").Value = "  +4.25%  "
$ws.Range("D13").Value = "'3.793.09"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "'27.83"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "'67.123.03"
$ws.Range("E16").Value = "  +3.92%  "
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").Value = "'3.236.55"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").Value = "'13.38"
$ws.Range("E20").Value = "  +2.81%  "
$ws.Range("D21").Value = "'372.68"
$ws.Range("E21").Value = "  +4.73%  "
$ws.Range("D22").Value = "'7.56"
$ws.Range("E22").Value = "  +4.06%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").Value = "'70.76"
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("D25").Value = "'0.509"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").Value = "'9.57"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("E28").Value = "  +2.78%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").Value = "'1.96"
$ws.Range("E30").Value = "  +3.51%  "
$ws.Range("D31").Value = "'5.64"
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "'1.26"
$ws.Range("E34").Value = "  +3.97%  "
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("D36").Value = "'162.35"
$ws.Range("E36").Value = "  +5.72%  "
$ws.Range("E37").Value = "  +3.34%  "
$ws.Range("D38").Value = "'0.854"
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("E39").Value = "  +6.66%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'26.63"
$ws.Range("E40").Value = "  +1.76%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'6.73"
$ws.Range("E41").Value = "  +11.50%  "
$ws.Range("D42").Value = "'2.61"
$ws.Range("E42").Value = "  +3.92%  "
$ws.Range("D43").Value = "'359.82"
$ws.Range("E43").Value = "  +12.73%  "
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("D45").Value = "'2.694.93"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").Value = "'25.43"
$ws.Range("E46").Value = "  +5.35%  "
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("E48").Value = "  +2.73%  "
$ws.Range("D49").Value = "'0.0278"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'0.992"
$ws.Range("E50").Value = "  +5.70%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.102"
$ws.Range("E51").Value = "  +0.29%  "
